$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report date range from 18-10-2025 to 19-10-2025
$ws.Range("I1").Value = "19-10-2025 00:00:00"
$ws.Range("K1").Value = "19-10-2025 00:00:00"

# Swap the data (columns B through G) between two rows, keeping the
# serial number in column A fixed for each row.
function Swap-RowData {
    param($ws, $row1, $row2)
    $cols = @("B","C","D","E","F","G")
    foreach ($c in $cols) {
        $v1 = $ws.Range("$c$row1").Value()
        $v2 = $ws.Range("$c$row2").Value()
        $ws.Range("$c$row1").Value = $v2
        $ws.Range("$c$row2").Value = $v1
    }
}

# Row pairs whose data got swapped between the two rows
$rowPairs = @(
    , @(279, 280)
    , @(313, 314)
    , @(316, 317)
    , @(372, 373)
    , @(379, 380)
    , @(419, 420)
    , @(421, 422)
    , @(431, 432)
    , @(457, 458)
    , @(536, 537)
    , @(601, 602)
    , @(687, 688)
    , @(720, 721)
    , @(859, 860)
    , @(889, 890)
)
foreach ($pair in $rowPairs) {
    Swap-RowData $ws $pair[0] $pair[1]
}

# Rows 350, 351, 352: a three-way rotation, where each row takes on
# the data (columns B through G) of the next row, wrapping around
# (350 <- 351, 351 <- 352, 352 <- 350).
$cols = @("B","C","D","E","F","G")
foreach ($c in $cols) {
    $v350 = $ws.Range("${c}350").Value()
    $v351 = $ws.Range("${c}351").Value()
    $v352 = $ws.Range("${c}352").Value()
    $ws.Range("${c}350").Value = $v351
    $ws.Range("${c}351").Value = $v352
    $ws.Range("${c}352").Value = $v350
}
